$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2025-08-13 23:38"
$ws.Range("B2").Value = "4aad170"
$ws.Range("C2").Value = "[SECURITY]: Complete RLS implementation with user context clients"
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 341
$ws.Range("F2").Value = 235
$ws.Range("G2").Value = "Complete RLS implementation - Critical security upgrade with user context clients"

$wb.Save()
